$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase the Schwierigkeitsgrad (difficulty weight) for the three
# existing "*" rows from 0.5 to 1.
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1

# Add two new documentation rows (15 and 16) for "OOD" entries.
$ws.Range("A15").Value = "OOD"
$ws.Range("B15").Value = "*"
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "Simon"

$ws.Range("A16").Value = "OOD"
$ws.Range("B16").Value = "*"
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "Niklas"

# Keep the active selection where the author last left it.
$ws.Range("C13").Select()

$wb.Save()
